$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# Convert state/municipality names from upper case to title case
$ws.Range("A2").Value = "Aguascalientes"
$ws.Range("B2").Value = "Aguascalientes"
$ws.Range("B3").Value = "Calvillo"
$ws.Range("B4").Value = "Total"
$ws.Range("A5").Value = "Campeche"
$ws.Range("B5").Value = "Escárcega"
$ws.Range("B6").Value = "Total"
$ws.Range("A7").Value = "Chiapas"
$ws.Range("B7").Value = "Acacoyagua"
$ws.Range("B8").Value = "Angel Albino Corzo"
$ws.Range("B9").Value = "El Bosque"
$ws.Range("B10").Value = "Frontera Comalapa"
$ws.Range("B11").Value = "Ixtapa"
$ws.Range("B12").Value = "Mapastepec"
$ws.Range("B13").Value = "Mitontic"
$ws.Range("B14").Value = "Ocosingo"
$ws.Range("B15").Value = "Ocotepec"
$ws.Range("B16").Value = "Palenque"
$ws.Range("B17").Value = "San Cristóbal De Las Casas"
$ws.Range("B18").Value = "Siltepec"
$ws.Range("B19").Value = "Simojovel"
$ws.Range("B20").Value = "Tapachula"
$ws.Range("B21").Value = "Tecpatán"
$ws.Range("B22").Value = "Tuxtla Gutiérrez"
$ws.Range("B23").Value = "Villaflores"
$ws.Range("B24").Value = "Total"
$ws.Range("A25").Value = "Chihuahua"
$ws.Range("B25").Value = "Buenaventura"
$ws.Range("B26").Value = "Chihuahua"
$ws.Range("B27").Value = "Cuauhtémoc"
$ws.Range("B28").Value = "Ignacio Zaragoza"
$ws.Range("B29").Value = "Juárez"
$ws.Range("B30").Value = "Nuevo Casas Grandes"
$ws.Range("B31").Value = "Saucillo"
$ws.Range("B32").Value = "Total"
$ws.Range("A33").Value = "Ciudad De México"
$ws.Range("B33").Value = "Benito Juárez"
$ws.Range("B34").Value = "Cuauhtémoc"
$ws.Range("B35").Value = "Gustavo A. Madero"
$ws.Range("B36").Value = "Iztapalapa"
$ws.Range("B37").Value = "Miguel Hidalgo"
$ws.Range("B38").Value = "Total"
$ws.Range("A39").Value = "Coahuila De Zaragoza"
$ws.Range("B39").Value = "Múzquiz"
$ws.Range("B40").Value = "Piedras Negras"
$ws.Range("B41").Value = "Torreón"
$ws.Range("B42").Value = "Total"
$ws.Range("A43").Value = "Durango"
$ws.Range("B43").Value = "Canatlán"
$ws.Range("B44").Value = "Durango"
$ws.Range("B45").Value = "Guadalupe Victoria"
$ws.Range("B46").Value = "Mapimí"
$ws.Range("B47").Value = "Vicente Guerrero"
$ws.Range("B48").Value = "Total"
$ws.Range("A49").Value = "Estado De México"
$ws.Range("B49").Value = "Acambay De Ruíz Castañeda"
$ws.Range("B50").Value = "Atlacomulco"
$ws.Range("B51").Value = "Ecatepec De Morelos"
$ws.Range("B52").Value = "Jocotitlán"
$ws.Range("B53").Value = "Nezahualcóyotl"
$ws.Range("B54").Value = "San Felipe Del Progreso"
$ws.Range("B55").Value = "Tenancingo"
$ws.Range("B56").Value = "Toluca"
$ws.Range("B57").Value = "Villa Guerrero"
$ws.Range("B58").Value = "Villa Victoria"
$ws.Range("B59").Value = "Total"
$ws.Range("A60").Value = "Guanajuato"
$ws.Range("B60").Value = "Acámbaro"
$ws.Range("B61").Value = "Apaseo El Alto"
$ws.Range("B62").Value = "Celaya"
$ws.Range("B63").Value = "Doctor Mora"
$ws.Range("B64").Value = "Guanajuato"
$ws.Range("B65").Value = "León"
$ws.Range("B66").Value = "Romita"
$ws.Range("B67").Value = "San Luis De La Paz"
$ws.Range("B68").Value = "San Miguel De Allende"
$ws.Range("B69").Value = "Santa Cruz De Juventino Rosas"
$ws.Range("B70").Value = "Silao De La Victoria"
$ws.Range("B71").Value = "Villagrán"
$ws.Range("B72").Value = "Total"
$ws.Range("A73").Value = "Guerrero"
$ws.Range("B73").Value = "Atoyac De Álvarez"
$ws.Range("B74").Value = "Chilpancingo De Los Bravo"
$ws.Range("B75").Value = "Coyuca De Benítez"
$ws.Range("B76").Value = "Tecoanapa"
$ws.Range("B77").Value = "Técpan De Galeana"
$ws.Range("B78").Value = "Zirándaro"
$ws.Range("B79").Value = "Total"
$ws.Range("A80").Value = "Hidalgo"
$ws.Range("B80").Value = "Ixmiquilpan"
$ws.Range("B81").Value = "Pachuca De Soto"
$ws.Range("B82").Value = "Tezontepec De Aldama"
$ws.Range("B83").Value = "Tulancingo De Bravo"
$ws.Range("B84").Value = "Total"
$ws.Range("A85").Value = "Jalisco"
$ws.Range("B85").Value = "Autlán De Navarro"
$ws.Range("B86").Value = "Ayutla"
$ws.Range("B87").Value = "Cabo Corrientes"
$ws.Range("B88").Value = "Casimiro Castillo"
$ws.Range("B89").Value = "Colotlán"
$ws.Range("B90").Value = "Guadalajara"
$ws.Range("B91").Value = "Jesús María"
$ws.Range("B92").Value = "Lagos De Moreno"
$ws.Range("B93").Value = "Ojuelos De Jalisco"
$ws.Range("B94").Value = "San Gabriel"
$ws.Range("B95").Value = "San Juan De Los Lagos"
$ws.Range("B96").Value = "Tequila"
$ws.Range("B97").Value = "Tizapán El Alto"
$ws.Range("B98").Value = "Tolimán"
$ws.Range("B99").Value = "Tomatlán"
$ws.Range("B100").Value = "Zapopan"
$ws.Range("B101").Value = "Total"
$ws.Range("A102").Value = "Michoacán De Ocampo"
$ws.Range("B102").Value = "Hidalgo"
$ws.Range("B103").Value = "Paracho"
$ws.Range("B104").Value = "Pátzcuaro"
$ws.Range("B105").Value = "Tlalpujahua"
$ws.Range("B106").Value = "Uruapan"
$ws.Range("B107").Value = "Total"
$ws.Range("A108").Value = "Nayarit"
$ws.Range("B108").Value = "Del Nayar"
$ws.Range("B109").Value = "San Blas"
$ws.Range("B110").Value = "Santiago Ixcuintla"
$ws.Range("B111").Value = "Tuxpan"
$ws.Range("B112").Value = "Total"
$ws.Range("A113").Value = "Nuevo León"
$ws.Range("B113").Value = "Sabinas Hidalgo"
$ws.Range("B114").Value = "Total"
$ws.Range("A115").Value = "Oaxaca"
$ws.Range("B115").Value = "San Juan Lachao"
$ws.Range("B116").Value = "San Lucas Zoquiápam"
$ws.Range("B117").Value = "Santa María Chimalapa"
$ws.Range("B118").Value = "Santa María Tonameca"
$ws.Range("B119").Value = "Tlacolula De Matamoros"
$ws.Range("B120").Value = "Total"
$ws.Range("A121").Value = "Puebla"
$ws.Range("B121").Value = "Ajalpan"
$ws.Range("B122").Value = "Chignahuapan"
$ws.Range("B123").Value = "Cuyoaco"
$ws.Range("B124").Value = "Jalpan"
$ws.Range("B125").Value = "Puebla"
$ws.Range("B126").Value = "San Martín Totoltepec"
$ws.Range("B127").Value = "Vicente Guerrero"
$ws.Range("B128").Value = "Zacatlán"
$ws.Range("B129").Value = "Total"
$ws.Range("A130").Value = "Querétaro"
$ws.Range("B130").Value = "Jalpan De Serra"
$ws.Range("B131").Value = "Querétaro"
$ws.Range("B132").Value = "Tequisquiapan"
$ws.Range("B133").Value = "Total"
$ws.Range("A134").Value = "San Luis Potosí"
$ws.Range("B134").Value = "Ciudad Valles"
$ws.Range("B135").Value = "Ebano"
$ws.Range("B136").Value = "Mexquitic De Carmona"
$ws.Range("B137").Value = "Salinas"
$ws.Range("B138").Value = "San Luis Potosí"
$ws.Range("B139").Value = "Villa De Ramos"
$ws.Range("B140").Value = "Xilitla"
$ws.Range("B141").Value = "Total"
$ws.Range("A142").Value = "Sinaloa"
$ws.Range("B142").Value = "Ahome"
$ws.Range("B143").Value = "Total"
$ws.Range("A144").Value = "Sonora"
$ws.Range("B144").Value = "Cajeme"
$ws.Range("B145").Value = "Total"
$ws.Range("A146").Value = "Tamaulipas"
$ws.Range("B146").Value = "Victoria"
$ws.Range("B147").Value = "Total"
$ws.Range("A148").Value = "Veracruz De Ignacio De La Llave"
$ws.Range("B148").Value = "Alpatláhuac"
$ws.Range("B149").Value = "Atlahuilco"
$ws.Range("B150").Value = "Calcahualco"
$ws.Range("B151").Value = "Coscomatepec"
$ws.Range("B152").Value = "Espinal"
$ws.Range("B153").Value = "Huayacocotla"
$ws.Range("B154").Value = "La Antigua"
$ws.Range("B155").Value = "Mixtla De Altamirano"
$ws.Range("B156").Value = "No Se Registró El Municipio/Condado/Alcaldía De Nacimiento"
$ws.Range("B157").Value = "Orizaba"
$ws.Range("B158").Value = "Veracruz"
$ws.Range("B159").Value = "Zongolica"
$ws.Range("B160").Value = "Total"
$ws.Range("A161").Value = "Zacatecas"
$ws.Range("B161").Value = "Fresnillo"
$ws.Range("B162").Value = "Guadalupe"
$ws.Range("B163").Value = "Ojocaliente"
$ws.Range("B164").Value = "Total"
$ws.Range("A165").Value = "Total"

# Remove footer/metadata rows (167-171)
$ws.Range("A167:A171").EntireRow.Delete()
